$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 250228.62
$ws.Range("I12").Value = 221.5
$ws.Range("K12").Value = 221.5
$ws.Range("M12").Value = -51.5

$ws.Range("H70").Value = 1461.7142
$ws.Range("I70").Value = 1398.4
$ws.Range("J70").Value = 1620
$ws.Range("K70").Value = 4195.200000000001
$ws.Range("L70").Value = 4860
$ws.Range("M70").Value = -3925.200000000001
$ws.Range("N70").Value = -5400

$ws.Range("H73").Value = 1461.7142
$ws.Range("I73").Value = 1398.4
$ws.Range("J73").Value = 1620
$ws.Range("K73").Value = 4195.200000000001
$ws.Range("L73").Value = 4860
$ws.Range("M73").Value = -3259.200000000001
$ws.Range("N73").Value = -6732

$ws.Range("H138").Value = 5692.413
$ws.Range("I138").Value = 2652.6072
$ws.Range("J138").Value = 7022.328
$ws.Range("K138").Value = 7957.821599999999
$ws.Range("L138").Value = 21066.984
$ws.Range("M138").Value = -2817.821599999999
$ws.Range("N138").Value = -31346.984

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1718.5862
$ws.Range("I74").Value = 1273.6
$ws.Range("K74").Value = 1273.6
$ws.Range("M74").Value = -399.5999999999999

$ws.Range("H77").Value = 1718.5862
$ws.Range("I77").Value = 1273.6
$ws.Range("K77").Value = 6368
$ws.Range("M77").Value = -2000

$ws.Range("H107").Value = 30000
$ws.Range("J107").Value = 30000
$ws.Range("L107").Value = 30000
$ws.Range("N107").Value = -37680

$ws.Range("H134").Value = 34900
$ws.Range("J134").Value = 34900
$ws.Range("L134").Value = 34900
$ws.Range("N134").Value = -45040

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1906.7667
$ws.Range("I105").Value = 1689.9546
$ws.Range("J105").Value = 2503
$ws.Range("K105").Value = 1689.9546
$ws.Range("L105").Value = 2503
$ws.Range("M105").Value = 57.04539999999997
$ws.Range("N105").Value = -5997

$ws.Range("H108").Value = 34000
$ws.Range("J108").Value = 34000
$ws.Range("L108").Value = 34000
$ws.Range("N108").Value = -41680

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1124
$ws.Range("I22").Value = 546.1818
$ws.Range("K22").Value = 546.1818
$ws.Range("M22").Value = -196.1818

$ws.Range("H62").Value = 3786.9092
$ws.Range("I62").Value = 2516.6667
$ws.Range("K62").Value = 2516.6667
$ws.Range("M62").Value = -1892.6667

$ws.Range("H65").Value = 3786.9092
$ws.Range("I65").Value = 2516.6667
$ws.Range("K65").Value = 12583.3335
$ws.Range("M65").Value = -9463.333500000001

$ws.Range("H105").Value = 3229.875
$ws.Range("I105").Value = 2778.5334
$ws.Range("J105").Value = 10000
$ws.Range("K105").Value = 2778.5334
$ws.Range("L105").Value = 10000
$ws.Range("M105").Value = -1031.5334
$ws.Range("N105").Value = -13494

$ws.Range("H108").Value = 31000
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 31000
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 31000
$ws.Range("M108").ClearContents()
$ws.Range("N108").Value = -38680

$ws.Range("H141").Value = 30802.941
$ws.Range("I141").Value = 7400
$ws.Range("J141").Value = 32265.625
$ws.Range("K141").Value = 7400
$ws.Range("L141").Value = 32265.625
$ws.Range("M141").Value = -2220
$ws.Range("N141").Value = -42625.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H61").Value = 8103.2666
$ws.Range("I61").Value = 104.5
$ws.Range("J61").Value = 9333.846
$ws.Range("K61").Value = 313.5
$ws.Range("L61").Value = 28001.538
$ws.Range("M61").Value = -98.5
$ws.Range("N61").Value = -28431.538

$ws.Range("H82").Value = 3046.8

$ws.Range("H85").Value = 3046.8

$ws.Range("H131").Value = 1461.807
$ws.Range("I131").Value = 5749.75
$ws.Range("J131").Value = 1138.1887
$ws.Range("K131").Value = 17249.25
$ws.Range("L131").Value = 3414.5661
$ws.Range("M131").Value = -12209.25
$ws.Range("N131").Value = -13494.5661

$ws.Range("H137").Value = 3143.2727
$ws.Range("I137").Value = 2579.375
$ws.Range("J137").Value = 4647
$ws.Range("K137").Value = 7738.125
$ws.Range("L137").Value = 13941
$ws.Range("M137").Value = -2638.125
$ws.Range("N137").Value = -24141

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3456.3333
$ws.Range("I80").Value = 3492.3076
$ws.Range("J80").Value = 3222.5
$ws.Range("K80").Value = 3492.3076
$ws.Range("L80").Value = 3222.5
$ws.Range("M80").Value = -2494.3076
$ws.Range("N80").Value = -5218.5

$ws.Range("H83").Value = 3456.3333
$ws.Range("I83").Value = 3492.3076
$ws.Range("J83").Value = 3222.5
$ws.Range("K83").Value = 17461.538
$ws.Range("L83").Value = 16112.5
$ws.Range("M83").Value = -12469.538
$ws.Range("N83").Value = -26096.5

$ws.Range("H126").Value = 912808.9399999999
$ws.Range("I126").Value = 1967.5
$ws.Range("J126").Value = 1433289.8
$ws.Range("K126").Value = 5902.5
$ws.Range("L126").Value = 4299869.4
$ws.Range("M126").Value = -3432.5
$ws.Range("N126").Value = -4304809.4

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1101.5883
$ws.Range("I22").Value = 390
$ws.Range("J22").Value = 1489.7273
$ws.Range("K22").Value = 390
$ws.Range("L22").Value = 1489.7273
$ws.Range("M22").Value = -95
$ws.Range("N22").Value = -2079.7273

$ws.Range("H27").Value = 1101.5883
$ws.Range("I27").Value = 390
$ws.Range("J27").Value = 1489.7273
$ws.Range("K27").Value = 390
$ws.Range("L27").Value = 1489.7273
$ws.Range("M27").Value = -283
$ws.Range("N27").Value = -1703.7273

$ws.Range("H46").Value = 1164.5853
$ws.Range("I46").Value = 919.1177
$ws.Range("J46").Value = 2356.8572
$ws.Range("K46").Value = 919.1177
$ws.Range("L46").Value = 2356.8572
$ws.Range("M46").Value = -731.1177
$ws.Range("N46").Value = -2732.8572

$ws.Range("H101").Value = 20590.5
$ws.Range("J101").Value = 20590.5
$ws.Range("L101").Value = 20590.5
$ws.Range("N101").Value = -27080.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 29589.578
$ws.Range("J103").Value = 29705.666
$ws.Range("L103").Value = 29705.666
$ws.Range("N103").Value = -32049.666

$ws.Range("H133").Value = 33357.5
$ws.Range("J133").Value = 33357.5
$ws.Range("L133").Value = 33357.5
$ws.Range("N133").Value = -43477.5

$ws.Range("H139").Value = 30000
$ws.Range("J139").Value = 30000
$ws.Range("L139").Value = 30000
$ws.Range("N139").Value = -40280
